# Applies the betexplorer scrape update (script run 26-11-2023 20:30):
#  - Row 58/59 (match order on 01/11/2023) are swapped.
#  - Rows 65/66/67 (matches on 09/11/2023) are cyclically rotated.
#  - Three new match rows (77, 78, 79) are appended for 23/11/2023 matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap/rotate the "match data" columns (F..V) between two rows,
# reading with Value2 (the only reliable read accessor in this runtime) and
# writing with Value (the reliable write accessor).
# ---------------------------------------------------------------------------
$cols = 6..22   # F=6 ... V=22

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value = $vals[$c]
    }
}

# --- Swap rows 58 and 59 (F:V) --------------------------------------------
$row58 = Get-RowValues 58
$row59 = Get-RowValues 59
Set-RowValues 58 $row59
Set-RowValues 59 $row58

# --- Rotate rows 65, 66, 67 (F:V): new65=old67, new66=old65, new67=old66 --
$row65 = Get-RowValues 65
$row66 = Get-RowValues 66
$row67 = Get-RowValues 67
Set-RowValues 65 $row67
Set-RowValues 66 $row65
Set-RowValues 67 $row66

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Append three new rows (77, 78, 79), copying the style of the last existing
# data row (76) so borders / fonts / number formats line up, then filling in
# the scraped values.
# ---------------------------------------------------------------------------
function Add-MatchRow($rowNum, $values) {
    # Clone formatting from row 76 (A:V) onto the new row first.
    $ws.Range("A76:V76").Copy()
    $ws.Range("A$rowNum`:V$rowNum").PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = $values[0]    # A indice
    $ws.Cells.Item($rowNum, 2).Value = $values[1]    # B pais
    $ws.Cells.Item($rowNum, 3).Value = $values[2]    # C torneio
    $ws.Cells.Item($rowNum, 4).Value = $values[3]    # D temporada
    $ws.Cells.Item($rowNum, 5).Value = $values[4]    # E data_partida
    $ws.Cells.Item($rowNum, 6).Value = $values[5]    # F home
    $ws.Cells.Item($rowNum, 7).Value = $values[6]    # G home_ft_gols
    $ws.Cells.Item($rowNum, 8).Value = $values[7]    # H away
    $ws.Cells.Item($rowNum, 9).Value = $values[8]    # I away_ft_gols
    $ws.Cells.Item($rowNum, 10).Value = $values[9]   # J home_opening_odds
    $ws.Cells.Item($rowNum, 11).Value = $values[10]  # K home_opening_data_hora
    $ws.Cells.Item($rowNum, 12).Value = $values[11]  # L home_closing_odds
    $ws.Cells.Item($rowNum, 13).Value = $values[12]  # M home_closing_data_hora
    $ws.Cells.Item($rowNum, 14).Value = $values[13]  # N draw_opening_odds
    $ws.Cells.Item($rowNum, 15).Value = $values[14]  # O draw_opening_data_hora
    $ws.Cells.Item($rowNum, 16).Value = $values[15]  # P draw_closing_odds
    $ws.Cells.Item($rowNum, 17).Value = $values[16]  # Q draw_closing_data_hora
    $ws.Cells.Item($rowNum, 18).Value = $values[17]  # R away_opening_odds
    $ws.Cells.Item($rowNum, 19).Value = $values[18]  # S away_opening_data_hora
    $ws.Cells.Item($rowNum, 20).Value = $values[19]  # T away_closing_odds
    $ws.Cells.Item($rowNum, 21).Value = $values[20]  # U away_closing_data_hora
    $ws.Cells.Item($rowNum, 22).Value = $values[21]  # V url_partida

    $excel.CutCopyMode = $false
}

Add-MatchRow 77 @(
    76, "iran", "persian-gulf-pro-league", "2023-2024", 45255.52083333334,
    "Tractor", 2, "Foolad", 0,
    1.74, "23/11/2023 09:13", 1.71, "25/11/2023 12:29",
    3.11, "23/11/2023 09:13", 2.63, "25/11/2023 12:29",
    4.63, "23/11/2023 09:13", 5.72, "25/11/2023 12:29",
    "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/tractor-foolad/fqonEWUf/"
)

Add-MatchRow 78 @(
    77, "iran", "persian-gulf-pro-league", "2023-2024", 45255.52083333334,
    "Shams Azar Qazvin", 1, "Gol Gohar", 1,
    2.52, "23/11/2023 09:13", 2.81, "25/11/2023 12:25",
    2.81, "23/11/2023 09:13", 2.54, "25/11/2023 12:25",
    2.76, "23/11/2023 09:13", 3.1, "25/11/2023 12:25",
    "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/shams-azar-qazvin-gol-gohar/0MorFCFl/"
)

Add-MatchRow 79 @(
    78, "iran", "persian-gulf-pro-league", "2023-2024", 45255.625,
    "Zob Ahan", 1, "Esteghlal F.C.", 1,
    3.32, "23/11/2023 09:13", 3.87, "25/11/2023 14:56",
    2.82, "23/11/2023 09:13", 2.6, "25/11/2023 14:59",
    2.18, "23/11/2023 09:13", 2.34, "25/11/2023 14:56",
    "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/zob-ahan-esteghlal-teh/8nkjDjp1/"
)
